$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "1.1.1.A" : insert a new "label" column (B) with footnoted fuel
# names, and append two extra rows (Bioenergy and waste / Total [note 6])
# that previously had no entry on this sheet.
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("1.1.1.A")
$wsA.Columns.Item(2).Insert()
$wsA.Cells.Item(1,2).Value = "label"
$wsA.Cells.Item(2,2).Value = "Coal [note 1]"
$wsA.Cells.Item(3,2).Value = "Petroleum [note 2]"
$wsA.Cells.Item(4,2).Value = "Natural gas [note 3]"
$wsA.Cells.Item(5,2).Value = "Nuclear electricity [note 4]"
$wsA.Cells.Item(6,2).Value = "Wind, solar and hydro electricity [note 4] [note 5]"
$wsA.Cells.Item(7,2).Value = "Net electricity imports"
$wsA.Cells.Item(7,3).Value = "Net electricity imports"
$wsA.Cells.Item(8,2).Value = "Bioenergy and waste"
$wsA.Cells.Item(9,2).Value = "Total [note 6]"

$wsA.Range("D7:D8").Select()

# ---------------------------------------------------------------------------
# Sheet "1.1.1.B" : same new "label" column, rows already existed.
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("1.1.1.B")
$wsB.Columns.Item(2).Insert()
$wsB.Cells.Item(1,2).Value = "label"
$wsB.Cells.Item(2,2).Value = "Coal [note 1]"
$wsB.Cells.Item(3,2).Value = "Petroleum [note 2]"
$wsB.Cells.Item(4,2).Value = "Natural gas [note 3]"
$wsB.Cells.Item(5,2).Value = "Nuclear electricity [note 4]"
$wsB.Cells.Item(6,2).Value = "Wind, solar and hydro electricity [note 4] [note 5]"
$wsB.Cells.Item(7,2).Value = "Net electricity imports"
$wsB.Cells.Item(7,3).Value = "Net electricity imports"
$wsB.Cells.Item(8,2).Value = "Bioenergy and waste"
$wsB.Cells.Item(8,3).Value = "Bioenergy and waste"
$wsB.Cells.Item(9,2).Value = "Total [note 6]"

$wsB.Range("B2:B9").Select()

# ---------------------------------------------------------------------------
# Sheet "1.1.1.C" : same new "label" column; also the custom column width
# that used to sit on column B now belongs to column C, and the
# "topLeftCell" scroll-freeze on the view goes away.
# ---------------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("1.1.1.C")
$wsC.Columns.Item(2).Insert()
$wsC.Cells.Item(1,2).Value = "label"
$wsC.Cells.Item(2,2).Value = "Coal [note 1]"
$wsC.Cells.Item(3,2).Value = "Petroleum [note 2]"
$wsC.Cells.Item(4,2).Value = "Natural gas [note 3]"
$wsC.Cells.Item(5,2).Value = "Nuclear electricity [note 4]"
$wsC.Cells.Item(6,2).Value = "Wind, solar and hydro electricity [note 4] [note 5]"
$wsC.Cells.Item(7,2).Value = "Net electricity imports"
$wsC.Cells.Item(7,3).Value = "Net electricity imports"
$wsC.Cells.Item(8,2).Value = "Bioenergy and waste"
$wsC.Cells.Item(8,3).Value = "Bioenergy and waste"
$wsC.Cells.Item(9,2).Value = "Total [note 6]"

$wsC.Range("B2").Select()

# ---------------------------------------------------------------------------
# View / selection bookkeeping: "1.1.6" becomes the active sheet/tab with a
# new zoom level and selection, so "1.2" (previously active) stops being
# the tab-selected sheet.
# ---------------------------------------------------------------------------
$ws16 = $wb.Worksheets.Item("1.1.6")
$ws16.Activate()
$ws16.Application.ActiveWindow.Zoom = 175
$ws16.Range("B5").Select()
